# Generate Report for Handback
#
# The handback transform for b1b319dc-36ea-4ae0-b081-7f493b12b11a.md failed
# (priority mismatch between the handback file path and the handoff type),
# while 14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.md is still "In Translation".
# Update the report rows for these two files accordingly across all three
# worksheets (Overview, zh-cn, de-de), and widen the Error Detail column so
# the new error message is readable.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$errMsgZhCn = "The handback priority in file path mt\b1b319dc-36ea-4ae0-b081-7f493b12b11a.a35dbeaf042692b6ab15e45f0535bafe8a451ee5.zh-cn.xlf is not match with handoff type ht."
$errMsgDeDe = "The handback priority in file path mt\b1b319dc-36ea-4ae0-b081-7f493b12b11a.a35dbeaf042692b6ab15e45f0535bafe8a451ee5.de-de.xlf is not match with handoff type ht."

# ----- Overview sheet: row 4 becomes b1b319dc, row 5 becomes 14ca0de6 -----
$ws1.Range("A4").Value = "b1b319dc-36ea-4ae0-b081-7f493b12b11a.md"
$ws1.Range("B4").Value = "e2e\b1b319dc-36ea-4ae0-b081-7f493b12b11a.md"
$ws1.Range("E4").Value = "Handback transform failed"
$ws1.Range("F4").Value = "Handback transform failed"
$ws1.Range("G4").Value = "2016-08-14 16:23:51"

$ws1.Range("A5").Value = "14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.md"
$ws1.Range("B5").Value = "e2e\14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.md"
$ws1.Range("E5").Value = "In Translation"
$ws1.Range("F5").Value = "In Translation"
$ws1.Range("G5").Value = "2016-08-14 16:19:03"

# ----- zh-cn sheet: row 4 becomes b1b319dc, row 5 becomes 14ca0de6 -----
$ws2.Range("A4").Value = "b1b319dc-36ea-4ae0-b081-7f493b12b11a.md"
$ws2.Range("C4").Value = "Handback transform failed"
$ws2.Range("G4").Value = "b1b319dc-36ea-4ae0-b081-7f493b12b11a.a35dbeaf042692b6ab15e45f0535bafe8a451ee5.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-14 16:23:43"
$ws2.Range("P4").Value = $errMsgZhCn

$ws2.Range("A5").Value = "14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.md"
$ws2.Range("C5").Value = "In Translation"
$ws2.Range("G5").Value = "14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.1a771575cbf5430b21b9b3cf98dae9a6097bca3c.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-14 16:18:53"

# ----- de-de sheet: row 4 becomes b1b319dc, row 5 becomes 14ca0de6 -----
$ws3.Range("A4").Value = "b1b319dc-36ea-4ae0-b081-7f493b12b11a.md"
$ws3.Range("C4").Value = "Handback transform failed"
$ws3.Range("G4").Value = "b1b319dc-36ea-4ae0-b081-7f493b12b11a.a35dbeaf042692b6ab15e45f0535bafe8a451ee5.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-14 16:23:51"
$ws3.Range("P4").Value = $errMsgDeDe

$ws3.Range("A5").Value = "14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.md"
$ws3.Range("C5").Value = "In Translation"
$ws3.Range("G5").Value = "14ca0de6-a1a0-4189-a75c-9ca025ed6c8d.1a771575cbf5430b21b9b3cf98dae9a6097bca3c.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-14 16:19:03"

# Widen the "Error Detail" column (P) on both language sheets so the new
# error message is visible.
$ws2.Columns.Item(16).ColumnWidth = 40
$ws3.Columns.Item(16).ColumnWidth = 40
